$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '64.473.42'
Set-TextValue $ws.Range("E2") '  +1.18%  '

Set-TextValue $ws.Range("D3") '2.630.46'
Set-TextValue $ws.Range("E3") '  +0.31%  '

Set-TextValue $ws.Range("E4") '  +0.05%  '

Set-TextValue $ws.Range("D5") '594.45'
Set-TextValue $ws.Range("E5") '  -0.24%  '

Set-TextValue $ws.Range("D6") '152.74'
Set-TextValue $ws.Range("E6") '  +1.82%  '

Set-TextValue $ws.Range("E7") '  +0.05%  '

Set-TextValue $ws.Range("D8") '0.589'
Set-TextValue $ws.Range("E8") '  -0.20%  '

Set-TextValue $ws.Range("E9") '  +4.26%  '

Set-TextValue $ws.Range("E10") '  +4.09%  '

Set-TextValue $ws.Range("D11") '5.78'
Set-TextValue $ws.Range("E11") '  +1.27%  '

Set-TextValue $ws.Range("E12") '  +1.07%  '

Set-TextValue $ws.Range("D13") '28.37'
Set-TextValue $ws.Range("E13") '  +2.28%  '

Set-TextValue $ws.Range("D14") '3.103.61'
Set-TextValue $ws.Range("E14") '  +0.30%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D15") '0.0000172'
Set-TextValue $ws.Range("E15") '  +13.01%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range("D16") '64.382.80'
Set-TextValue $ws.Range("E16") '  +1.24%  '

Set-TextValue $ws.Range("D17") '2.653.24'
Set-TextValue $ws.Range("E17") '  +0.89%  '

Set-TextValue $ws.Range("D18") '12.29'
Set-TextValue $ws.Range("E18") '  -0.37%  '

Set-TextValue $ws.Range("D19") '4.78'
Set-TextValue $ws.Range("E19") '  +2.54%  '

Set-TextValue $ws.Range("D20") '349.36'
Set-TextValue $ws.Range("E20") '  +0.68%  '

Set-TextValue $ws.Range("D21") '7.10'
Set-TextValue $ws.Range("E21") '  +3.58%  '

Set-TextValue $ws.Range("E22") '  +0.26%  '

Set-TextValue $ws.Range("D23") '67.39'
Set-TextValue $ws.Range("E23") '  +1.51%  '

Set-TextValue $ws.Range("D24") '1.70'
Set-TextValue $ws.Range("E24") '  -0.71%  '

Set-TextValue $ws.Range("D25") '9.27'
Set-TextValue $ws.Range("E25") '  +0.46%  '

Set-TextValue $ws.Range("D26") '1.66'
Set-TextValue $ws.Range("E26") '  -0.32%  '

Set-TextValue $ws.Range("D27") '8.24'
Set-TextValue $ws.Range("E27") '  +0.09%  '

Set-TextValue $ws.Range("D28") '551.21'
Set-TextValue $ws.Range("E28") '  +0.42%  '

Set-TextValue $ws.Range("E29") '  +0.40%  '

Set-TextValue $ws.Range("D30") '1.00'
Set-TextValue $ws.Range("E30") '  +0.07%  '

Set-TextValue $ws.Range("D31") '0.0₃0908'
Set-TextValue $ws.Range("E31") '  +7.67%  '

Set-TextValue $ws.Range("E32") '  +2.00%  '

Set-TextValue $ws.Range("D33") '1.79'
Set-TextValue $ws.Range("E33") '  +2.96%  '

Set-TextValue $ws.Range("D34") '5.54'
Set-TextValue $ws.Range("E34") '  +5.77%  '

Set-TextValue $ws.Range("D35") '6.20'
Set-TextValue $ws.Range("E35") '  +1.45%  '

Set-TextValue $ws.Range("D36") '0.422'
Set-TextValue $ws.Range("E36") '  +2.78%  '

Set-TextValue $ws.Range("D37") '164.54'
Set-TextValue $ws.Range("E37") '  -2.41%  '

$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D38") '20.15'
Set-TextValue $ws.Range("E38") '  +3.85%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D39") '2.01'
Set-TextValue $ws.Range("E39") '  +3.51%  '

Set-TextValue $ws.Range("E40") '  -0.07%  '

Set-TextValue $ws.Range("E41") '  -0.05%  '

Set-TextValue $ws.Range("D42") '168.11'
Set-TextValue $ws.Range("E42") '  +0.75%  '

Set-TextValue $ws.Range("D43") '41.46'
Set-TextValue $ws.Range("E43") '  +4.04%  '

Set-TextValue $ws.Range("D44") '4.10'
Set-TextValue $ws.Range("E44") '  +4.47%  '

Set-TextValue $ws.Range("D45") '23.38'
Set-TextValue $ws.Range("E45") '  +8.97%  '

$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D46") '0.0590'
Set-TextValue $ws.Range("E46") '  -0.33%  '

$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D47") '2.21'
Set-TextValue $ws.Range("E47") '  +12.15%  '

Set-TextValue $ws.Range("D48") '0.642'
Set-TextValue $ws.Range("E48") '  +1.80%  '

Set-TextValue $ws.Range("D49") '0.0251'
Set-TextValue $ws.Range("E49") '  +1.07%  '

Set-TextValue $ws.Range("D50") '0.0978'
Set-TextValue $ws.Range("E50") '  +1.21%  '

Set-TextValue $ws.Range("E51") '  +0.17%  '
